$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.644.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.06%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.617.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.58%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '628.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.55'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.79%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.617.24'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.52%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.497'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.39%  '

$ws.Range("E10").Value = '  +9.35%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.35'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.441'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000229'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.49'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.66%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.228.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.57%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '70.079.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.610.63'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.77%  '

$ws.Range("E18").Value = '  +0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.06'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.69%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +14.93%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '463.92'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.57%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.647'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.78%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.71'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.20%  '

$ws.Range("E25").Value = '  +14.42%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.761.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.39%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +14.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.64'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.34%  '

$ws.Range("E31").Value = '  +11.04%  '

$ws.Range("E32").Value = '  +14.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.57'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.51%  '

$ws.Range("E35").Value = '  +4.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.88%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.614.67'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.79%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.48'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.77%  '

$ws.Range("E39").Value = '  +13.01%  '

$ws.Range("E40").Value = '  +0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '180.04'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0928'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.30%  '

$ws.Range("E44").Value = '  +4.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '32.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +21.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.917'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.77%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +15.82%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.06'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.48%  '

$ws.Range("E49").Value = '  +14.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.83'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.01%  '

$ws.Range("E51").Value = '  +10.27%  '

